# Updates cryptocurrency price (column D) and volume/change (column E) values
# on Sheet1, rows 2-51, with the latest scrape from coinranking.com.
#
# A leading apostrophe is prepended to every literal so that Excel stores the
# cell as literal text (quote-prefixed), matching the original inlineStr cell
# type, instead of auto-converting numeric-looking strings (e.g. "216.34" or
# "1.668.16") into actual numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.047.52"
$ws.Range('E2').Value = "'  -1.92%  "
$ws.Range('D3').Value = "'1.665.03"
$ws.Range('E3').Value = "'  -1.29%  "
$ws.Range('D5').Value = "'216.34"
$ws.Range('E5').Value = "'  -0.46%  "
$ws.Range('D6').Value = "'0.5101"
$ws.Range('E6').Value = "'  +1.92%  "
$ws.Range('E7').Value = "'  +0.21%  "
$ws.Range('D8').Value = "'0.2626"
$ws.Range('E8').Value = "'  +0.31%  "
$ws.Range('D9').Value = "'0.06404"
$ws.Range('E9').Value = "'  +3.14%  "
$ws.Range('E10').Value = "'  -1.26%  "
$ws.Range('D11').Value = "'0.07421"
$ws.Range('E11').Value = "'  +1.86%  "
$ws.Range('D12').Value = "'1.671.31"
$ws.Range('E13').Value = "'  +1.15%  "
$ws.Range('D14').Value = "'0.5799"
$ws.Range('E14').Value = "'  +0.54%  "
$ws.Range('D15').Value = "'0.000008507"
$ws.Range('E15').Value = "'  +3.81%  "
$ws.Range('D16').Value = "'64.26"
$ws.Range('E16').Value = "'  -0.79%  "
$ws.Range('D17').Value = "'26.114.55"
$ws.Range('E17').Value = "'  -1.73%  "
$ws.Range('D18').Value = "'4.892"
$ws.Range('E18').Value = "'  -2.16%  "
$ws.Range('E20').Value = "'  -0.36%  "
$ws.Range('D21').Value = "'188.41"
$ws.Range('E21').Value = "'  +1.63%  "
$ws.Range('D22').Value = "'6.195"
$ws.Range('E22').Value = "'  -0.02%  "
$ws.Range('E23').Value = "'  +0.18%  "
$ws.Range('D24').Value = "'145.56"
$ws.Range('E24').Value = "'  +0.68%  "
$ws.Range('D25').Value = "'7.604"
$ws.Range('E25').Value = "'  +1.41%  "
$ws.Range('E26').Value = "'  +4.65%  "
$ws.Range('E27').Value = "'  +0.70%  "
$ws.Range('D28').Value = "'0.06484"
$ws.Range('E28').Value = "'  +14.02%  "
$ws.Range('D29').Value = "'1.308"
$ws.Range('E29').Value = "'  +0.68%  "
$ws.Range('D30').Value = "'1.314"
$ws.Range('E30').Value = "'  -0.41%  "
$ws.Range('D31').Value = "'3.524"
$ws.Range('D32').Value = "'3.504"
$ws.Range('E32').Value = "'  +0.69%  "
$ws.Range('D33').Value = "'1.627"
$ws.Range('E33').Value = "'  -0.70%  "
$ws.Range('D34').Value = "'1.016"
$ws.Range('E34').Value = "'  +0.71%  "
$ws.Range('D35').Value = "'0.6047"
$ws.Range('E35').Value = "'  +2.04%  "
$ws.Range('E36').Value = "'  -0.15%  "
$ws.Range('E37').Value = "'  +1.85%  "
$ws.Range('D38').Value = "'6.196"
$ws.Range('E38').Value = "'  +5.21%  "
$ws.Range('D39').Value = "'0.01608"
$ws.Range('E39').Value = "'  +1.08%  "
$ws.Range('D40').Value = "'1.074.38"
$ws.Range('E40').Value = "'  +0.40%  "
$ws.Range('D41').Value = "'0.8579"
$ws.Range('E41').Value = "'  +0.12%  "
$ws.Range('E42').Value = "'  +0.80%  "
$ws.Range('D43').Value = "'100.49"
$ws.Range('E43').Value = "'  +2.42%  "
$ws.Range('D44').Value = "'1.812.97"
$ws.Range('E44').Value = "'  -1.61%  "
$ws.Range('D45').Value = "'0.00000000115"
$ws.Range('E45').Value = "'  +8.30%  "
$ws.Range('D46').Value = "'56.08"
$ws.Range('E46').Value = "'  -0.52%  "
$ws.Range('D47').Value = "'1.003"
$ws.Range('E47').Value = "'  -0.15%  "
$ws.Range('D48').Value = "'8.007"
$ws.Range('E48').Value = "'  -0.07%  "
$ws.Range('E49').Value = "'  +0.19%  "
$ws.Range('D50').Value = "'0.4295"
$ws.Range('E50').Value = "'  -0.39%  "
$ws.Range('D51').Value = "'5.929"
$ws.Range('E51').Value = "'  +4.59%  "
